$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44893
$ws.Cells.Item(2, 9).Value = 'Primera'
$ws.Cells.Item(2, 10).Value = 120
$ws.Cells.Item(2, 11).Value = 9000
$ws.Cells.Item(2, 12).Value = 9500
$ws.Cells.Item(2, 13).Value = 9250
$ws.Cells.Item(2, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(2, 16).Value = 370

$ws.Cells.Item(3, 4).Value = 44159
$ws.Cells.Item(3, 9).Value = 'Primera'
$ws.Cells.Item(3, 10).Value = 42
$ws.Cells.Item(3, 11).Value = 6500
$ws.Cells.Item(3, 12).Value = 7000
$ws.Cells.Item(3, 13).Value = 6738
$ws.Cells.Item(3, 15).Value = 'Región del Maule'
$ws.Cells.Item(3, 16).Value = 270

$ws.Cells.Item(4, 4).Value = 44509
$ws.Cells.Item(4, 9).Value = 'Primera'
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 8000
$ws.Cells.Item(4, 12).Value = 9000
$ws.Cells.Item(4, 13).Value = 8500
$ws.Cells.Item(4, 15).Value = 'Región del Maule'
$ws.Cells.Item(4, 16).Value = 340

$ws.Cells.Item(5, 4).Value = 44809
$ws.Cells.Item(5, 9).Value = 'Primera'
$ws.Cells.Item(5, 10).Value = 60
$ws.Cells.Item(5, 11).Value = 12000
$ws.Cells.Item(5, 12).Value = 13000
$ws.Cells.Item(5, 13).Value = 12500
$ws.Cells.Item(5, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(5, 16).Value = 500

$ws.Cells.Item(6, 4).Value = 44495
$ws.Cells.Item(6, 9).Value = 'Primera'
$ws.Cells.Item(6, 10).Value = 60
$ws.Cells.Item(6, 11).Value = 8000
$ws.Cells.Item(6, 12).Value = 9000
$ws.Cells.Item(6, 13).Value = 8500
$ws.Cells.Item(6, 15).Value = 'Región del Maule'
$ws.Cells.Item(6, 16).Value = 340

$ws.Cells.Item(7, 4).Value = 44803
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 50
$ws.Cells.Item(7, 11).Value = 13000
$ws.Cells.Item(7, 12).Value = 14000
$ws.Cells.Item(7, 13).Value = 13500
$ws.Cells.Item(7, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(7, 16).Value = 540

$ws.Cells.Item(8, 4).Value = 44529
$ws.Cells.Item(8, 9).Value = 'Primera'
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 6000
$ws.Cells.Item(8, 12).Value = 7000
$ws.Cells.Item(8, 13).Value = 6500
$ws.Cells.Item(8, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(8, 16).Value = 260

$ws.Cells.Item(9, 4).Value = 44553
$ws.Cells.Item(9, 9).Value = 'Primera'
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 6500
$ws.Cells.Item(9, 12).Value = 7000
$ws.Cells.Item(9, 13).Value = 6750
$ws.Cells.Item(9, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(9, 16).Value = 270

$ws.Cells.Item(10, 4).Value = 44530
$ws.Cells.Item(10, 9).Value = 'Primera'
$ws.Cells.Item(10, 10).Value = 100
$ws.Cells.Item(10, 11).Value = 6000
$ws.Cells.Item(10, 12).Value = 7000
$ws.Cells.Item(10, 13).Value = 6500
$ws.Cells.Item(10, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(10, 16).Value = 260

$ws.Cells.Item(11, 4).Value = 44845
$ws.Cells.Item(11, 9).Value = 'Primera'
$ws.Cells.Item(11, 10).Value = 120
$ws.Cells.Item(11, 11).Value = 9000
$ws.Cells.Item(11, 12).Value = 9500
$ws.Cells.Item(11, 13).Value = 9250
$ws.Cells.Item(11, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(11, 16).Value = 370

$ws.Cells.Item(12, 4).Value = 44512
$ws.Cells.Item(12, 9).Value = 'Primera'
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 7000
$ws.Cells.Item(12, 12).Value = 8000
$ws.Cells.Item(12, 13).Value = 7500
$ws.Cells.Item(12, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(12, 16).Value = 300

$ws.Cells.Item(13, 4).Value = 45258
$ws.Cells.Item(13, 9).Value = 'Primera'
$ws.Cells.Item(13, 10).Value = 100
$ws.Cells.Item(13, 11).Value = 10000
$ws.Cells.Item(13, 12).Value = 10000
$ws.Cells.Item(13, 13).Value = 10000
$ws.Cells.Item(13, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(13, 16).Value = 400

$ws.Cells.Item(14, 4).Value = 44918
$ws.Cells.Item(14, 9).Value = 'Primera'
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 15000
$ws.Cells.Item(14, 12).Value = 16000
$ws.Cells.Item(14, 13).Value = 15500
$ws.Cells.Item(14, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(14, 16).Value = 620

$ws.Cells.Item(15, 4).Value = 45225
$ws.Cells.Item(15, 9).Value = 'Primera'
$ws.Cells.Item(15, 10).Value = 30
$ws.Cells.Item(15, 11).Value = 12000
$ws.Cells.Item(15, 12).Value = 12000
$ws.Cells.Item(15, 13).Value = 12000
$ws.Cells.Item(15, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(15, 16).Value = 480

$ws.Cells.Item(16, 4).Value = 44167
$ws.Cells.Item(16, 9).Value = 'Primera'
$ws.Cells.Item(16, 10).Value = 60
$ws.Cells.Item(16, 11).Value = 8000
$ws.Cells.Item(16, 12).Value = 9000
$ws.Cells.Item(16, 13).Value = 8500
$ws.Cells.Item(16, 15).Value = 'Región del Maule'
$ws.Cells.Item(16, 16).Value = 340

$ws.Cells.Item(17, 4).Value = 45218
$ws.Cells.Item(17, 9).Value = 'Primera'
$ws.Cells.Item(17, 10).Value = 80
$ws.Cells.Item(17, 11).Value = 15000
$ws.Cells.Item(17, 12).Value = 15000
$ws.Cells.Item(17, 13).Value = 15000
$ws.Cells.Item(17, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(17, 16).Value = 600

$ws.Cells.Item(18, 4).Value = 44799
$ws.Cells.Item(18, 9).Value = 'Primera'
$ws.Cells.Item(18, 10).Value = 30
$ws.Cells.Item(18, 11).Value = 11000
$ws.Cells.Item(18, 12).Value = 11000
$ws.Cells.Item(18, 13).Value = 11000
$ws.Cells.Item(18, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(18, 16).Value = 440

$ws.Cells.Item(19, 4).Value = 44813
$ws.Cells.Item(19, 9).Value = 'Primera'
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 11).Value = 12000
$ws.Cells.Item(19, 12).Value = 13000
$ws.Cells.Item(19, 13).Value = 12500
$ws.Cells.Item(19, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(19, 16).Value = 500

$ws.Cells.Item(20, 4).Value = 44855
$ws.Cells.Item(20, 9).Value = 'Primera'
$ws.Cells.Item(20, 10).Value = 120
$ws.Cells.Item(20, 11).Value = 7000
$ws.Cells.Item(20, 12).Value = 7500
$ws.Cells.Item(20, 13).Value = 7250
$ws.Cells.Item(20, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(20, 16).Value = 290

$ws.Cells.Item(21, 4).Value = 44517
$ws.Cells.Item(21, 9).Value = 'Primera'
$ws.Cells.Item(21, 10).Value = 100
$ws.Cells.Item(21, 11).Value = 6000
$ws.Cells.Item(21, 12).Value = 7000
$ws.Cells.Item(21, 13).Value = 6500
$ws.Cells.Item(21, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(21, 16).Value = 260

$ws.Cells.Item(22, 4).Value = 44516
$ws.Cells.Item(22, 9).Value = 'Primera'
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 7000
$ws.Cells.Item(22, 12).Value = 8000
$ws.Cells.Item(22, 13).Value = 7500
$ws.Cells.Item(22, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(22, 16).Value = 300

$ws.Cells.Item(23, 4).Value = 45140
$ws.Cells.Item(23, 9).Value = 'Primera'
$ws.Cells.Item(23, 10).Value = 30
$ws.Cells.Item(23, 11).Value = 15000
$ws.Cells.Item(23, 12).Value = 15000
$ws.Cells.Item(23, 13).Value = 15000
$ws.Cells.Item(23, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(23, 16).Value = 600

$ws.Cells.Item(24, 4).Value = 45211
$ws.Cells.Item(24, 9).Value = 'Primera'
$ws.Cells.Item(24, 10).Value = 60
$ws.Cells.Item(24, 11).Value = 12000
$ws.Cells.Item(24, 12).Value = 12000
$ws.Cells.Item(24, 13).Value = 12000
$ws.Cells.Item(24, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(24, 16).Value = 480

$ws.Cells.Item(25, 4).Value = 44900
$ws.Cells.Item(25, 9).Value = 'Primera'
$ws.Cells.Item(25, 10).Value = 80
$ws.Cells.Item(25, 11).Value = 12000
$ws.Cells.Item(25, 12).Value = 12000
$ws.Cells.Item(25, 13).Value = 12000
$ws.Cells.Item(25, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(25, 16).Value = 480

$ws.Cells.Item(26, 4).Value = 44524
$ws.Cells.Item(26, 9).Value = 'Primera'
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 6000
$ws.Cells.Item(26, 12).Value = 7000
$ws.Cells.Item(26, 13).Value = 6500
$ws.Cells.Item(26, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(26, 16).Value = 260

$ws.Cells.Item(27, 4).Value = 44819
$ws.Cells.Item(27, 9).Value = 'Primera'
$ws.Cells.Item(27, 10).Value = 80
$ws.Cells.Item(27, 11).Value = 12000
$ws.Cells.Item(27, 12).Value = 13000
$ws.Cells.Item(27, 13).Value = 12500
$ws.Cells.Item(27, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(27, 16).Value = 500

$ws.Cells.Item(28, 4).Value = 45189
$ws.Cells.Item(28, 9).Value = 'Primera'
$ws.Cells.Item(28, 10).Value = 60
$ws.Cells.Item(28, 11).Value = 15000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 13).Value = 15000
$ws.Cells.Item(28, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(28, 16).Value = 600

$ws.Cells.Item(29, 4).Value = 45217
$ws.Cells.Item(29, 9).Value = 'Primera'
$ws.Cells.Item(29, 10).Value = 100
$ws.Cells.Item(29, 11).Value = 15000
$ws.Cells.Item(29, 12).Value = 15000
$ws.Cells.Item(29, 13).Value = 15000
$ws.Cells.Item(29, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(29, 16).Value = 600

$ws.Cells.Item(30, 4).Value = 45260
$ws.Cells.Item(30, 9).Value = 'Primera'
$ws.Cells.Item(30, 10).Value = 100
$ws.Cells.Item(30, 11).Value = 10000
$ws.Cells.Item(30, 12).Value = 10000
$ws.Cells.Item(30, 13).Value = 10000
$ws.Cells.Item(30, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(30, 16).Value = 400

$ws.Cells.Item(31, 4).Value = 45240
$ws.Cells.Item(31, 9).Value = 'Primera'
$ws.Cells.Item(31, 10).Value = 60
$ws.Cells.Item(31, 11).Value = 10000
$ws.Cells.Item(31, 12).Value = 10000
$ws.Cells.Item(31, 13).Value = 10000
$ws.Cells.Item(31, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(31, 16).Value = 400

$ws.Cells.Item(32, 4).Value = 44161
$ws.Cells.Item(32, 9).Value = 'Primera'
$ws.Cells.Item(32, 10).Value = 53
$ws.Cells.Item(32, 11).Value = 6500
$ws.Cells.Item(32, 12).Value = 7000
$ws.Cells.Item(32, 13).Value = 6764
$ws.Cells.Item(32, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(32, 16).Value = 271

$ws.Cells.Item(33, 4).Value = 44537
$ws.Cells.Item(33, 9).Value = 'Primera'
$ws.Cells.Item(33, 10).Value = 60
$ws.Cells.Item(33, 11).Value = 6500
$ws.Cells.Item(33, 12).Value = 7000
$ws.Cells.Item(33, 13).Value = 6750
$ws.Cells.Item(33, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(33, 16).Value = 270

$ws.Cells.Item(34, 4).Value = 44448
$ws.Cells.Item(34, 9).Value = 'Primera'
$ws.Cells.Item(34, 10).Value = 60
$ws.Cells.Item(34, 11).Value = 14000
$ws.Cells.Item(34, 12).Value = 15000
$ws.Cells.Item(34, 13).Value = 14500
$ws.Cells.Item(34, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(34, 16).Value = 580

$ws.Cells.Item(35, 4).Value = 45204
$ws.Cells.Item(35, 9).Value = 'Primera'
$ws.Cells.Item(35, 10).Value = 50
$ws.Cells.Item(35, 11).Value = 15000
$ws.Cells.Item(35, 12).Value = 15000
$ws.Cells.Item(35, 13).Value = 15000
$ws.Cells.Item(35, 15).Value = 'Región del Maule'
$ws.Cells.Item(35, 16).Value = 600

$ws.Cells.Item(36, 4).Value = 44523
$ws.Cells.Item(36, 9).Value = 'Primera'
$ws.Cells.Item(36, 10).Value = 80
$ws.Cells.Item(36, 11).Value = 6000
$ws.Cells.Item(36, 12).Value = 7000
$ws.Cells.Item(36, 13).Value = 6500
$ws.Cells.Item(36, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(36, 16).Value = 260

$ws.Cells.Item(37, 4).Value = 44536
$ws.Cells.Item(37, 9).Value = 'Primera'
$ws.Cells.Item(37, 10).Value = 80
$ws.Cells.Item(37, 11).Value = 6500
$ws.Cells.Item(37, 12).Value = 7000
$ws.Cells.Item(37, 13).Value = 6750
$ws.Cells.Item(37, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(37, 16).Value = 270

$ws.Cells.Item(38, 4).Value = 44873
$ws.Cells.Item(38, 9).Value = 'Primera'
$ws.Cells.Item(38, 10).Value = 80
$ws.Cells.Item(38, 11).Value = 6500
$ws.Cells.Item(38, 12).Value = 7000
$ws.Cells.Item(38, 13).Value = 6750
$ws.Cells.Item(38, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(38, 16).Value = 270

$ws.Cells.Item(39, 4).Value = 45131
$ws.Cells.Item(39, 9).Value = 'Primera'
$ws.Cells.Item(39, 10).Value = 40
$ws.Cells.Item(39, 11).Value = 16000
$ws.Cells.Item(39, 12).Value = 16000
$ws.Cells.Item(39, 13).Value = 16000
$ws.Cells.Item(39, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(39, 16).Value = 640

$ws.Cells.Item(40, 4).Value = 45222
$ws.Cells.Item(40, 9).Value = 'Primera'
$ws.Cells.Item(40, 10).Value = 60
$ws.Cells.Item(40, 11).Value = 15000
$ws.Cells.Item(40, 12).Value = 15000
$ws.Cells.Item(40, 13).Value = 15000
$ws.Cells.Item(40, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(40, 16).Value = 600

$ws.Cells.Item(41, 4).Value = 45091
$ws.Cells.Item(41, 9).Value = 'Primera'
$ws.Cells.Item(41, 10).Value = 60
$ws.Cells.Item(41, 11).Value = 15000
$ws.Cells.Item(41, 12).Value = 15000
$ws.Cells.Item(41, 13).Value = 15000
$ws.Cells.Item(41, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(41, 16).Value = 600

$ws.Cells.Item(42, 4).Value = 45247
$ws.Cells.Item(42, 9).Value = 'Primera'
$ws.Cells.Item(42, 10).Value = 200
$ws.Cells.Item(42, 11).Value = 10000
$ws.Cells.Item(42, 12).Value = 10000
$ws.Cells.Item(42, 13).Value = 10000
$ws.Cells.Item(42, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(42, 16).Value = 400

$ws.Cells.Item(43, 4).Value = 44882
$ws.Cells.Item(43, 9).Value = 'Primera'
$ws.Cells.Item(43, 10).Value = 120
$ws.Cells.Item(43, 11).Value = 7000
$ws.Cells.Item(43, 12).Value = 7500
$ws.Cells.Item(43, 13).Value = 7250
$ws.Cells.Item(43, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(43, 16).Value = 290

$ws.Cells.Item(44, 4).Value = 44910
$ws.Cells.Item(44, 9).Value = 'Primera'
$ws.Cells.Item(44, 10).Value = 100
$ws.Cells.Item(44, 11).Value = 10000
$ws.Cells.Item(44, 12).Value = 12000
$ws.Cells.Item(44, 13).Value = 11000
$ws.Cells.Item(44, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(44, 16).Value = 440

$ws.Cells.Item(45, 4).Value = 45246
$ws.Cells.Item(45, 9).Value = 'Primera'
$ws.Cells.Item(45, 10).Value = 100
$ws.Cells.Item(45, 11).Value = 10000
$ws.Cells.Item(45, 12).Value = 10000
$ws.Cells.Item(45, 13).Value = 10000
$ws.Cells.Item(45, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(45, 16).Value = 400

$ws.Cells.Item(46, 4).Value = 44806
$ws.Cells.Item(46, 9).Value = 'Primera'
$ws.Cells.Item(46, 10).Value = 60
$ws.Cells.Item(46, 11).Value = 13000
$ws.Cells.Item(46, 12).Value = 14000
$ws.Cells.Item(46, 13).Value = 13500
$ws.Cells.Item(46, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(46, 16).Value = 540

$ws.Cells.Item(47, 4).Value = 45203
$ws.Cells.Item(47, 9).Value = 'Primera'
$ws.Cells.Item(47, 10).Value = 120
$ws.Cells.Item(47, 11).Value = 14000
$ws.Cells.Item(47, 12).Value = 15000
$ws.Cells.Item(47, 13).Value = 14667
$ws.Cells.Item(47, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(47, 16).Value = 587

$ws.Cells.Item(48, 4).Value = 44484
$ws.Cells.Item(48, 9).Value = 'Primera'
$ws.Cells.Item(48, 10).Value = 30
$ws.Cells.Item(48, 11).Value = 8500
$ws.Cells.Item(48, 12).Value = 9000
$ws.Cells.Item(48, 13).Value = 8750
$ws.Cells.Item(48, 15).Value = 'Región del Maule'
$ws.Cells.Item(48, 16).Value = 350

$ws.Cells.Item(49, 4).Value = 45251
$ws.Cells.Item(49, 9).Value = 'Primera'
$ws.Cells.Item(49, 10).Value = 120
$ws.Cells.Item(49, 11).Value = 9500
$ws.Cells.Item(49, 12).Value = 10000
$ws.Cells.Item(49, 13).Value = 9750
$ws.Cells.Item(49, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(49, 16).Value = 390

$ws.Cells.Item(50, 4).Value = 44519
$ws.Cells.Item(50, 9).Value = 'Primera'
$ws.Cells.Item(50, 10).Value = 80
$ws.Cells.Item(50, 11).Value = 6000
$ws.Cells.Item(50, 12).Value = 7000
$ws.Cells.Item(50, 13).Value = 6500
$ws.Cells.Item(50, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(50, 16).Value = 260

$ws.Cells.Item(51, 4).Value = 44812
$ws.Cells.Item(51, 9).Value = 'Primera'
$ws.Cells.Item(51, 10).Value = 60
$ws.Cells.Item(51, 11).Value = 12000
$ws.Cells.Item(51, 12).Value = 13000
$ws.Cells.Item(51, 13).Value = 12500
$ws.Cells.Item(51, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(51, 16).Value = 500

$ws.Cells.Item(52, 4).Value = 44491
$ws.Cells.Item(52, 9).Value = 'Primera'
$ws.Cells.Item(52, 10).Value = 60
$ws.Cells.Item(52, 11).Value = 8000
$ws.Cells.Item(52, 12).Value = 9000
$ws.Cells.Item(52, 13).Value = 8500
$ws.Cells.Item(52, 15).Value = 'Región del Maule'
$ws.Cells.Item(52, 16).Value = 340

$ws.Cells.Item(53, 4).Value = 45239
$ws.Cells.Item(53, 9).Value = 'Primera'
$ws.Cells.Item(53, 10).Value = 80
$ws.Cells.Item(53, 11).Value = 10000
$ws.Cells.Item(53, 12).Value = 10000
$ws.Cells.Item(53, 13).Value = 10000
$ws.Cells.Item(53, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(53, 16).Value = 400

$ws.Cells.Item(54, 4).Value = 44540
$ws.Cells.Item(54, 9).Value = 'Primera'
$ws.Cells.Item(54, 10).Value = 100
$ws.Cells.Item(54, 11).Value = 6500
$ws.Cells.Item(54, 12).Value = 7000
$ws.Cells.Item(54, 13).Value = 6750
$ws.Cells.Item(54, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(54, 16).Value = 270

$ws.Cells.Item(55, 4).Value = 45141
$ws.Cells.Item(55, 9).Value = 'Primera'
$ws.Cells.Item(55, 10).Value = 50
$ws.Cells.Item(55, 11).Value = 15000
$ws.Cells.Item(55, 12).Value = 15000
$ws.Cells.Item(55, 13).Value = 15000
$ws.Cells.Item(55, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(55, 16).Value = 600

$ws.Cells.Item(56, 4).Value = 44522
$ws.Cells.Item(56, 9).Value = 'Primera'
$ws.Cells.Item(56, 10).Value = 100
$ws.Cells.Item(56, 11).Value = 6000
$ws.Cells.Item(56, 12).Value = 7000
$ws.Cells.Item(56, 13).Value = 6500
$ws.Cells.Item(56, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(56, 16).Value = 260

$ws.Cells.Item(57, 4).Value = 44488
$ws.Cells.Item(57, 9).Value = 'Primera'
$ws.Cells.Item(57, 10).Value = 60
$ws.Cells.Item(57, 11).Value = 8000
$ws.Cells.Item(57, 12).Value = 9000
$ws.Cells.Item(57, 13).Value = 8500
$ws.Cells.Item(57, 15).Value = 'Región del Maule'
$ws.Cells.Item(57, 16).Value = 340

$ws.Cells.Item(58, 4).Value = 44482
$ws.Cells.Item(58, 9).Value = 'Primera'
$ws.Cells.Item(58, 10).Value = 120
$ws.Cells.Item(58, 11).Value = 8000
$ws.Cells.Item(58, 12).Value = 9000
$ws.Cells.Item(58, 13).Value = 8500
$ws.Cells.Item(58, 15).Value = 'Región del Maule'
$ws.Cells.Item(58, 16).Value = 340

$ws.Cells.Item(59, 4).Value = 45079
$ws.Cells.Item(59, 9).Value = 'Primera'
$ws.Cells.Item(59, 10).Value = 30
$ws.Cells.Item(59, 11).Value = 15000
$ws.Cells.Item(59, 12).Value = 15000
$ws.Cells.Item(59, 13).Value = 15000
$ws.Cells.Item(59, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(59, 16).Value = 600

$ws.Cells.Item(60, 4).Value = 44511
$ws.Cells.Item(60, 9).Value = 'Primera'
$ws.Cells.Item(60, 10).Value = 100
$ws.Cells.Item(60, 11).Value = 7000
$ws.Cells.Item(60, 12).Value = 8000
$ws.Cells.Item(60, 13).Value = 7500
$ws.Cells.Item(60, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(60, 16).Value = 300

$ws.Cells.Item(61, 4).Value = 44876
$ws.Cells.Item(61, 9).Value = 'Primera'
$ws.Cells.Item(61, 10).Value = 100
$ws.Cells.Item(61, 11).Value = 6500
$ws.Cells.Item(61, 12).Value = 7000
$ws.Cells.Item(61, 13).Value = 6750
$ws.Cells.Item(61, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(61, 16).Value = 270

$ws.Cells.Item(62, 4).Value = 44476
$ws.Cells.Item(62, 9).Value = 'Primera'
$ws.Cells.Item(62, 10).Value = 160
$ws.Cells.Item(62, 11).Value = 7500
$ws.Cells.Item(62, 12).Value = 8000
$ws.Cells.Item(62, 13).Value = 7750
$ws.Cells.Item(62, 15).Value = 'Región del Maule'
$ws.Cells.Item(62, 16).Value = 310

$ws.Cells.Item(63, 4).Value = 44166
$ws.Cells.Item(63, 9).Value = 'Primera'
$ws.Cells.Item(63, 10).Value = 56
$ws.Cells.Item(63, 11).Value = 7500
$ws.Cells.Item(63, 12).Value = 8000
$ws.Cells.Item(63, 13).Value = 7804
$ws.Cells.Item(63, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(63, 16).Value = 312

$ws.Cells.Item(64, 4).Value = 44533
$ws.Cells.Item(64, 9).Value = 'Primera'
$ws.Cells.Item(64, 10).Value = 80
$ws.Cells.Item(64, 11).Value = 6500
$ws.Cells.Item(64, 12).Value = 7000
$ws.Cells.Item(64, 13).Value = 6750
$ws.Cells.Item(64, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(64, 16).Value = 270

$ws.Cells.Item(65, 4).Value = 44487
$ws.Cells.Item(65, 9).Value = 'Primera'
$ws.Cells.Item(65, 10).Value = 30
$ws.Cells.Item(65, 11).Value = 8000
$ws.Cells.Item(65, 12).Value = 8000
$ws.Cells.Item(65, 13).Value = 8000
$ws.Cells.Item(65, 15).Value = 'Región del Maule'
$ws.Cells.Item(65, 16).Value = 320

$ws.Cells.Item(66, 4).Value = 44487
$ws.Cells.Item(66, 9).Value = 'Segunda'
$ws.Cells.Item(66, 10).Value = 30
$ws.Cells.Item(66, 11).Value = 9000
$ws.Cells.Item(66, 12).Value = 9000
$ws.Cells.Item(66, 13).Value = 9000
$ws.Cells.Item(66, 15).Value = 'Región del Maule'
$ws.Cells.Item(66, 16).Value = 360

$ws.Cells.Item(67, 4).Value = 44847
$ws.Cells.Item(67, 9).Value = 'Primera'
$ws.Cells.Item(67, 10).Value = 80
$ws.Cells.Item(67, 11).Value = 9000
$ws.Cells.Item(67, 12).Value = 9000
$ws.Cells.Item(67, 13).Value = 9000
$ws.Cells.Item(67, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(67, 16).Value = 360

$ws.Cells.Item(68, 4).Value = 45142
$ws.Cells.Item(68, 9).Value = 'Primera'
$ws.Cells.Item(68, 10).Value = 30
$ws.Cells.Item(68, 11).Value = 15000
$ws.Cells.Item(68, 12).Value = 15000
$ws.Cells.Item(68, 13).Value = 15000
$ws.Cells.Item(68, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(68, 16).Value = 600

$ws.Cells.Item(69, 4).Value = 45194
$ws.Cells.Item(69, 9).Value = 'Primera'
$ws.Cells.Item(69, 10).Value = 60
$ws.Cells.Item(69, 11).Value = 14000
$ws.Cells.Item(69, 12).Value = 14000
$ws.Cells.Item(69, 13).Value = 14000
$ws.Cells.Item(69, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(69, 16).Value = 560

$ws.Cells.Item(70, 4).Value = 44566
$ws.Cells.Item(70, 9).Value = 'Primera'
$ws.Cells.Item(70, 10).Value = 60
$ws.Cells.Item(70, 11).Value = 7000
$ws.Cells.Item(70, 12).Value = 7500
$ws.Cells.Item(70, 13).Value = 7250
$ws.Cells.Item(70, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(70, 16).Value = 290

$ws.Cells.Item(71, 4).Value = 44908
$ws.Cells.Item(71, 9).Value = 'Primera'
$ws.Cells.Item(71, 10).Value = 80
$ws.Cells.Item(71, 11).Value = 12000
$ws.Cells.Item(71, 12).Value = 12000
$ws.Cells.Item(71, 13).Value = 12000
$ws.Cells.Item(71, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(71, 16).Value = 480

$ws.Cells.Item(72, 4).Value = 45243
$ws.Cells.Item(72, 9).Value = 'Primera'
$ws.Cells.Item(72, 10).Value = 100
$ws.Cells.Item(72, 11).Value = 10000
$ws.Cells.Item(72, 12).Value = 10000
$ws.Cells.Item(72, 13).Value = 10000
$ws.Cells.Item(72, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(72, 16).Value = 400

$ws.Cells.Item(73, 4).Value = 44816
$ws.Cells.Item(73, 9).Value = 'Primera'
$ws.Cells.Item(73, 10).Value = 80
$ws.Cells.Item(73, 11).Value = 12000
$ws.Cells.Item(73, 12).Value = 13000
$ws.Cells.Item(73, 13).Value = 12500
$ws.Cells.Item(73, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(73, 16).Value = 500

$ws.Cells.Item(74, 4).Value = 44504
$ws.Cells.Item(74, 9).Value = 'Primera'
$ws.Cells.Item(74, 10).Value = 60
$ws.Cells.Item(74, 11).Value = 8000
$ws.Cells.Item(74, 12).Value = 9000
$ws.Cells.Item(74, 13).Value = 8500
$ws.Cells.Item(74, 15).Value = 'Región del Maule'
$ws.Cells.Item(74, 16).Value = 340

$ws.Cells.Item(75, 4).Value = 44859
$ws.Cells.Item(75, 9).Value = 'Primera'
$ws.Cells.Item(75, 10).Value = 120
$ws.Cells.Item(75, 11).Value = 6000
$ws.Cells.Item(75, 12).Value = 6500
$ws.Cells.Item(75, 13).Value = 6250
$ws.Cells.Item(75, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(75, 16).Value = 250

$ws.Cells.Item(76, 4).Value = 44518
$ws.Cells.Item(76, 9).Value = 'Primera'
$ws.Cells.Item(76, 10).Value = 60
$ws.Cells.Item(76, 11).Value = 6000
$ws.Cells.Item(76, 12).Value = 7000
$ws.Cells.Item(76, 13).Value = 6500
$ws.Cells.Item(76, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(76, 16).Value = 260

$ws.Cells.Item(77, 4).Value = 45264
$ws.Cells.Item(77, 9).Value = 'Primera'
$ws.Cells.Item(77, 10).Value = 60
$ws.Cells.Item(77, 11).Value = 10000
$ws.Cells.Item(77, 12).Value = 10000
$ws.Cells.Item(77, 13).Value = 10000
$ws.Cells.Item(77, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(77, 16).Value = 400

$ws.Cells.Item(78, 4).Value = 44924
$ws.Cells.Item(78, 9).Value = 'Primera'
$ws.Cells.Item(78, 10).Value = 60
$ws.Cells.Item(78, 11).Value = 15000
$ws.Cells.Item(78, 12).Value = 15000
$ws.Cells.Item(78, 13).Value = 15000
$ws.Cells.Item(78, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(78, 16).Value = 600

$ws.Cells.Item(79, 4).Value = 44875
$ws.Cells.Item(79, 9).Value = 'Primera'
$ws.Cells.Item(79, 10).Value = 120
$ws.Cells.Item(79, 11).Value = 6500
$ws.Cells.Item(79, 12).Value = 7000
$ws.Cells.Item(79, 13).Value = 6750
$ws.Cells.Item(79, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(79, 16).Value = 270

$ws.Cells.Item(80, 4).Value = 44860
$ws.Cells.Item(80, 9).Value = 'Primera'
$ws.Cells.Item(80, 10).Value = 120
$ws.Cells.Item(80, 11).Value = 6000
$ws.Cells.Item(80, 12).Value = 6500
$ws.Cells.Item(80, 13).Value = 6250
$ws.Cells.Item(80, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(80, 16).Value = 250

$ws.Cells.Item(81, 4).Value = 44515
$ws.Cells.Item(81, 9).Value = 'Primera'
$ws.Cells.Item(81, 10).Value = 100
$ws.Cells.Item(81, 11).Value = 7000
$ws.Cells.Item(81, 12).Value = 8000
$ws.Cells.Item(81, 13).Value = 7500
$ws.Cells.Item(81, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(81, 16).Value = 300

$ws.Cells.Item(82, 4).Value = 44473
$ws.Cells.Item(82, 9).Value = 'Primera'
$ws.Cells.Item(82, 10).Value = 60
$ws.Cells.Item(82, 11).Value = 9500
$ws.Cells.Item(82, 12).Value = 10000
$ws.Cells.Item(82, 13).Value = 9750
$ws.Cells.Item(82, 15).Value = 'Región del Maule'
$ws.Cells.Item(82, 16).Value = 390

$ws.Cells.Item(83, 4).Value = 45224
$ws.Cells.Item(83, 9).Value = 'Primera'
$ws.Cells.Item(83, 10).Value = 60
$ws.Cells.Item(83, 11).Value = 12000
$ws.Cells.Item(83, 12).Value = 13000
$ws.Cells.Item(83, 13).Value = 12500
$ws.Cells.Item(83, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(83, 16).Value = 500

$ws.Cells.Item(84, 4).Value = 44162
$ws.Cells.Item(84, 9).Value = 'Primera'
$ws.Cells.Item(84, 10).Value = 80
$ws.Cells.Item(84, 11).Value = 7000
$ws.Cells.Item(84, 12).Value = 8000
$ws.Cells.Item(84, 13).Value = 7562
$ws.Cells.Item(84, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(84, 16).Value = 302

$ws.Cells.Item(85, 4).Value = 44789
$ws.Cells.Item(85, 9).Value = 'Primera'
$ws.Cells.Item(85, 10).Value = 60
$ws.Cells.Item(85, 11).Value = 11000
$ws.Cells.Item(85, 12).Value = 12000
$ws.Cells.Item(85, 13).Value = 11500
$ws.Cells.Item(85, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(85, 16).Value = 460

$ws.Cells.Item(86, 4).Value = 44895
$ws.Cells.Item(86, 9).Value = 'Primera'
$ws.Cells.Item(86, 10).Value = 60
$ws.Cells.Item(86, 11).Value = 12000
$ws.Cells.Item(86, 12).Value = 12000
$ws.Cells.Item(86, 13).Value = 12000
$ws.Cells.Item(86, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(86, 16).Value = 480

$ws.Cells.Item(87, 4).Value = 44526
$ws.Cells.Item(87, 9).Value = 'Primera'
$ws.Cells.Item(87, 10).Value = 100
$ws.Cells.Item(87, 11).Value = 6000
$ws.Cells.Item(87, 12).Value = 7000
$ws.Cells.Item(87, 13).Value = 6500
$ws.Cells.Item(87, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(87, 16).Value = 260

$ws.Cells.Item(88, 4).Value = 44879
$ws.Cells.Item(88, 9).Value = 'Primera'
$ws.Cells.Item(88, 10).Value = 120
$ws.Cells.Item(88, 11).Value = 6500
$ws.Cells.Item(88, 12).Value = 7000
$ws.Cells.Item(88, 13).Value = 6750
$ws.Cells.Item(88, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(88, 16).Value = 270

$ws.Cells.Item(89, 4).Value = 44922
$ws.Cells.Item(89, 9).Value = 'Primera'
$ws.Cells.Item(89, 10).Value = 80
$ws.Cells.Item(89, 11).Value = 15000
$ws.Cells.Item(89, 12).Value = 15000
$ws.Cells.Item(89, 13).Value = 15000
$ws.Cells.Item(89, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(89, 16).Value = 600

$ws.Cells.Item(90, 4).Value = 44160
$ws.Cells.Item(90, 9).Value = 'Primera'
$ws.Cells.Item(90, 10).Value = 80
$ws.Cells.Item(90, 11).Value = 6500
$ws.Cells.Item(90, 12).Value = 7000
$ws.Cells.Item(90, 13).Value = 6688
$ws.Cells.Item(90, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(90, 16).Value = 268

$ws.Cells.Item(91, 4).Value = 44830
$ws.Cells.Item(91, 9).Value = 'Primera'
$ws.Cells.Item(91, 10).Value = 100
$ws.Cells.Item(91, 11).Value = 9000
$ws.Cells.Item(91, 12).Value = 9500
$ws.Cells.Item(91, 13).Value = 9250
$ws.Cells.Item(91, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(91, 16).Value = 370

$ws.Cells.Item(92, 4).Value = 45236
$ws.Cells.Item(92, 9).Value = 'Primera'
$ws.Cells.Item(92, 10).Value = 60
$ws.Cells.Item(92, 11).Value = 10000
$ws.Cells.Item(92, 12).Value = 10000
$ws.Cells.Item(92, 13).Value = 10000
$ws.Cells.Item(92, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(92, 16).Value = 400

$ws.Cells.Item(93, 4).Value = 45138
$ws.Cells.Item(93, 9).Value = 'Primera'
$ws.Cells.Item(93, 10).Value = 30
$ws.Cells.Item(93, 11).Value = 15000
$ws.Cells.Item(93, 12).Value = 15000
$ws.Cells.Item(93, 13).Value = 15000
$ws.Cells.Item(93, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(93, 16).Value = 600

$ws.Cells.Item(94, 4).Value = 45132
$ws.Cells.Item(94, 9).Value = 'Primera'
$ws.Cells.Item(94, 10).Value = 25
$ws.Cells.Item(94, 11).Value = 16000
$ws.Cells.Item(94, 12).Value = 16000
$ws.Cells.Item(94, 13).Value = 16000
$ws.Cells.Item(94, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(94, 16).Value = 640

$ws.Cells.Item(95, 4).Value = 44489
$ws.Cells.Item(95, 9).Value = 'Primera'
$ws.Cells.Item(95, 10).Value = 60
$ws.Cells.Item(95, 11).Value = 8000
$ws.Cells.Item(95, 12).Value = 9000
$ws.Cells.Item(95, 13).Value = 8500
$ws.Cells.Item(95, 15).Value = 'Región del Maule'
$ws.Cells.Item(95, 16).Value = 340

$ws.Cells.Item(96, 4).Value = 44165
$ws.Cells.Item(96, 9).Value = 'Primera'
$ws.Cells.Item(96, 10).Value = 38
$ws.Cells.Item(96, 11).Value = 8000
$ws.Cells.Item(96, 12).Value = 8500
$ws.Cells.Item(96, 13).Value = 8263
$ws.Cells.Item(96, 15).Value = 'Región del Maule'
$ws.Cells.Item(96, 16).Value = 331

$ws.Cells.Item(97, 4).Value = 44466
$ws.Cells.Item(97, 9).Value = 'Primera'
$ws.Cells.Item(97, 10).Value = 60
$ws.Cells.Item(97, 11).Value = 11000
$ws.Cells.Item(97, 12).Value = 12000
$ws.Cells.Item(97, 13).Value = 11500
$ws.Cells.Item(97, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(97, 16).Value = 460

$ws.Cells.Item(98, 4).Value = 44837
$ws.Cells.Item(98, 9).Value = 'Primera'
$ws.Cells.Item(98, 10).Value = 30
$ws.Cells.Item(98, 11).Value = 9000
$ws.Cells.Item(98, 12).Value = 9000
$ws.Cells.Item(98, 13).Value = 9000
$ws.Cells.Item(98, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(98, 16).Value = 360

$ws.Cells.Item(99, 4).Value = 44837
$ws.Cells.Item(99, 9).Value = 'Segunda'
$ws.Cells.Item(99, 10).Value = 30
$ws.Cells.Item(99, 11).Value = 9500
$ws.Cells.Item(99, 12).Value = 9500
$ws.Cells.Item(99, 13).Value = 9500
$ws.Cells.Item(99, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(99, 16).Value = 380

$ws.Cells.Item(100, 4).Value = 45209
$ws.Cells.Item(100, 9).Value = 'Primera'
$ws.Cells.Item(100, 10).Value = 50
$ws.Cells.Item(100, 11).Value = 14000
$ws.Cells.Item(100, 12).Value = 14000
$ws.Cells.Item(100, 13).Value = 14000
$ws.Cells.Item(100, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(100, 16).Value = 560

$ws.Cells.Item(101, 4).Value = 44883
$ws.Cells.Item(101, 9).Value = 'Primera'
$ws.Cells.Item(101, 10).Value = 120
$ws.Cells.Item(101, 11).Value = 7000
$ws.Cells.Item(101, 12).Value = 7500
$ws.Cells.Item(101, 13).Value = 7250
$ws.Cells.Item(101, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(101, 16).Value = 290

$ws.Cells.Item(102, 4).Value = 44797
$ws.Cells.Item(102, 9).Value = 'Primera'
$ws.Cells.Item(102, 10).Value = 60
$ws.Cells.Item(102, 11).Value = 12000
$ws.Cells.Item(102, 12).Value = 12000
$ws.Cells.Item(102, 13).Value = 12000
$ws.Cells.Item(102, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(102, 16).Value = 480

$ws.Cells.Item(103, 4).Value = 44792
$ws.Cells.Item(103, 9).Value = 'Primera'
$ws.Cells.Item(103, 10).Value = 50
$ws.Cells.Item(103, 11).Value = 12000
$ws.Cells.Item(103, 12).Value = 12000
$ws.Cells.Item(103, 13).Value = 12000
$ws.Cells.Item(103, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(103, 16).Value = 480

$ws.Cells.Item(104, 4).Value = 44897
$ws.Cells.Item(104, 9).Value = 'Primera'
$ws.Cells.Item(104, 10).Value = 120
$ws.Cells.Item(104, 11).Value = 12000
$ws.Cells.Item(104, 12).Value = 12500
$ws.Cells.Item(104, 13).Value = 12250
$ws.Cells.Item(104, 15).Value = 'Provincia de Diguillín'
$ws.Cells.Item(104, 16).Value = 490
